# Append the new transaction row (row 43) to Sheet1, mirroring the
# "Added row for JAHNAVI KOLASANI" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$r = 43

# Numeric columns -----------------------------------------------------
$ws.Cells.Item($r, 1).Value  = 15599            # admission_no
$ws.Cells.Item($r, 3).Value  = 9347559040       # phone
$ws.Cells.Item($r, 4).Value  = 11250            # amount
$ws.Cells.Item($r, 12).Value = 100000036600     # customer_acc_no
$ws.Cells.Item($r, 14).Value = 753702           # merchant_id
$ws.Cells.Item($r, 15).Value = 1234             # client_code
$ws.Cells.Item($r, 16).Value = 11000316561861   # atom_txn_id
$ws.Cells.Item($r, 17).Value = 1763044851       # merchant_txn_id
$ws.Cells.Item($r, 18).Value = 108566739973     # bank_ref_no
$ws.Cells.Item($r, 29).Value = 0                # sb_cess
$ws.Cells.Item($r, 30).Value = 0                # krishi_kalyan_cess
$ws.Cells.Item($r, 35).Value = 19295            # udfex1

# Plain text columns ----------------------------------------------------
$ws.Cells.Item($r, 2).Value  = "JAHNAVI KOLASANI"            # student_name
$ws.Cells.Item($r, 6).Value  = "13-Nov-2025 20:24:29"        # txn_date
$ws.Cells.Item($r, 8).Value  = "TRANSACTION IS SUCCESSFUL"   # description
$ws.Cells.Item($r, 9).Value  = "OK"                          # txn_status
$ws.Cells.Item($r, 10).Value = "Multi"                       # product
$ws.Cells.Item($r, 11).Value = "eleven thousand two hundred fifty" # amount_in_rupees
$ws.Cells.Item($r, 13).Value = "SALESIAN EDUCATION SOCIETY"  # merchant_name
$ws.Cells.Item($r, 19).Value = "INR"                         # currency
$ws.Cells.Item($r, 20).Value = "sale"                        # txn_type
$ws.Cells.Item($r, 21).Value = "ICICI UPI QR"                # bank_name
$ws.Cells.Item($r, 22).Value = "NRNS"                        # recon_status
$ws.Cells.Item($r, 23).Value = "IFSC0000000"                 # ifsc_code
$ws.Cells.Item($r, 24).Value = "MERCHANT"                    # merchant_type
$ws.Cells.Item($r, 25).Value = "UPI"                         # discriminator
$ws.Cells.Item($r, 26).Value = "kotakschoolvsp@gmail.com"    # email
$ws.Cells.Item($r, 34).Value = "REGULAR"                     # settlement_type

# Text columns that would otherwise be misread as numbers (comma-joined
# id lists) -- force Text format first so the comma isn't treated as a
# thousands separator.
$ws.Cells.Item($r, 36).NumberFormat = "@"
$ws.Cells.Item($r, 36).Value = "263081,264584"               # udfex2
$ws.Cells.Item($r, 37).NumberFormat = "@"
$ws.Cells.Item($r, 37).Value = "2037,2044"                   # udfex3

# Columns that are blank ("" / empty inline string) in the source row:
# net_amount_to_be_paid (E), settlement_date (G), txn_charges (AA),
# gst_18 (AB), total_chargeable (AE), beneficiary_name (AF),
# imps_status (AG), qr_transaction_type (AL) are intentionally left
# unset (blank cell), matching the empty values in the source data.
